$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# --- TextBox 4 ("Enterprise Grid Security" / "http://www.gaards.org") ----
# Collapse the two paragraphs into a single paragraph "Enterprise Grid Security"
# and bump the second half ("Security") up to the same 28pt size as the title.
$titleShape = $s1.Shapes.Item(2)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Text = "Enterprise Grid Security"
$titleRange.Characters(17, 8).Font.Size = 28

# Reposition/resize the title textbox (shrinks now that it is one line).
$titleShape.Left = 3276600 / 12700.0
$titleShape.Top = 152400 / 12700.0

# --- TextBox 5 ("Stephen Langella, ... Joel Saltz.") is removed ----------
$s1.Shapes.Item(3).Delete()

# --- New textbox: "http://www.cagrid.org" (right-aligned, top-right) -----
$cagridBox = $s1.Shapes.AddTextbox(1, 5181600 / 12700.0, 987623 / 12700.0, 2590800 / 12700.0, 307777 / 12700.0)
$cagridBox.Name = "TextBox 9"
$cagridBox.TextFrame.WordWrap = -1
$cagridBox.TextFrame.AutoSize = 1
$cagridBox.Fill.Visible = 0
$cagridRange = $cagridBox.TextFrame.TextRange
$cagridRange.Text = "http://www.cagrid.org"
$cagridRange.ParagraphFormat.Alignment = 3
$cagridRange.Font.Size = 14
$cagridRange.Font.Bold = -1
$cagridRange.Font.Color.RGB = 192
$cagridRange.Font.Name = "Calibri"

# --- New textbox: "GAARDS Security Infrastructure" (top-left) ------------
$gaardsBox = $s1.Shapes.AddTextbox(1, -76200 / 12700.0, 987623 / 12700.0, 2590800 / 12700.0, 307777 / 12700.0)
$gaardsBox.Name = "TextBox 10"
$gaardsBox.TextFrame.WordWrap = -1
$gaardsBox.TextFrame.AutoSize = 1
$gaardsBox.Fill.Visible = 0
$gaardsRange = $gaardsBox.TextFrame.TextRange
$gaardsRange.Text = "GAARDS Security Infrastructure"
$gaardsRange.Font.Size = 14
$gaardsRange.Font.Bold = -1
$gaardsRange.Font.Color.RGB = 192
$gaardsRange.Font.Name = "Calibri"

# ---------------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$bodyShape = $s2.Shapes.Item(3)
$bodyRange = $bodyShape.TextFrame.TextRange
$full = $bodyRange.Text
$idx = $full.IndexOf("Question ") + 1
$bodyRange.Characters($idx, 12).Text = "Question or "
# Text edits on this auto-fit shape recompute its height; restore the
# original (the authored change did not alter this shape's size).
$bodyShape.Height = 672.0
